$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The pagemap used backslash-separated "folder\file.php" names for the
# User-Visible Page column (B). Subfolders are no longer needed, so these
# page names get their backslash replaced by an underscore, e.g.
# "search\users.php" -> "search_users.php". (The includes\*.php and
# include_utils\*.php rows keep their real relative paths, so they are
# left untouched.)
$rowsToFix = @(12, 13, 14, 15, 16, 17, 18, 19, 20, 24, 25, 26, 28, 29, 30, 31, 33, 35, 37, 38, 39, 40, 41, 42)
foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val.Replace("\", "_")
    }
}

# Update the saved view state (scroll position + active selection) to match.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B49").Select()
